$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.129.21'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.558.20'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  +1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.780.68'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.562.69'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.094.59'
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.67'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.77'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0700'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.71'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.99'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.14'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.433.49'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.530'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.90'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.804'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.02'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.692.64'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.39'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0992'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0949'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.94%  '
